{"js": "// Replace each two-digit-divided-by-one-digit expression in the practice\n// table with the new expression dictated by the target diff. Every old\n// value is unique in the document, so a simple matchCase search + full\n// replace (\"Replace\" location) for each pair is unambiguous.\nconst replacements = [\n  [\"99\u00f79=\", \"10\u00f77=\"],\n  [\"51\u00f73=\", \"13\u00f77=\"],\n  [\"55\u00f74=\", \"99\u00f75=\"],\n  [\"96\u00f74=\", \"11\u00f79=\"],\n  [\"90\u00f75=\", \"79\u00f72=\"],\n  [\"85\u00f73=\", \"33\u00f79=\"],\n  [\"33\u00f73=\", \"18\u00f78=\"],\n  [\"48\u00f77=\", \"26\u00f75=\"],\n  [\"77\u00f72=\", \"63\u00f74=\"],\n  [\"36\u00f73=\", \"71\u00f79=\"],\n  [\"16\u00f76=\", \"51\u00f75=\"],\n  [\"34\u00f78=\", \"52\u00f75=\"],\n  [\"76\u00f78=\", \"16\u00f77=\"],\n  [\"81\u00f77=\", \"71\u00f73=\"],\n  [\"18\u00f79=\", \"56\u00f75=\"],\n  [\"11\u00f77=\", \"13\u00f79=\"],\n  [\"19\u00f79=\", \"94\u00f76=\"],\n  [\"35\u00f76=\", \"27\u00f72=\"],\n  [\"93\u00f75=\", \"61\u00f77=\"],\n  [\"35\u00f78=\", \"30\u00f74=\"],\n  [\"21\u00f75=\", \"32\u00f73=\"],\n  [\"29\u00f77=\", \"56\u00f79=\"],\n  [\"80\u00f79=\", \"89\u00f78=\"],\n  [\"40\u00f78=\", \"26\u00f79=\"],\n  [\"18\u00f77=\", \"60\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit-divided-by-one-digit expression in the practice\n# table with the new expression dictated by the target diff. Every old\n# value is unique in the document, so Find/Replace (wdReplaceAll, one hit\n# each) for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"99\u00f79=\", \"10\u00f77=\"),\n    @(\"51\u00f73=\", \"13\u00f77=\"),\n    @(\"55\u00f74=\", \"99\u00f75=\"),\n    @(\"96\u00f74=\", \"11\u00f79=\"),\n    @(\"90\u00f75=\", \"79\u00f72=\"),\n    @(\"85\u00f73=\", \"33\u00f79=\"),\n    @(\"33\u00f73=\", \"18\u00f78=\"),\n    @(\"48\u00f77=\", \"26\u00f75=\"),\n    @(\"77\u00f72=\", \"63\u00f74=\"),\n    @(\"36\u00f73=\", \"71\u00f79=\"),\n    @(\"16\u00f76=\", \"51\u00f75=\"),\n    @(\"34\u00f78=\", \"52\u00f75=\"),\n    @(\"76\u00f78=\", \"16\u00f77=\"),\n    @(\"81\u00f77=\", \"71\u00f73=\"),\n    @(\"18\u00f79=\", \"56\u00f75=\"),\n    @(\"11\u00f77=\", \"13\u00f79=\"),\n    @(\"19\u00f79=\", \"94\u00f76=\"),\n    @(\"35\u00f76=\", \"27\u00f72=\"),\n    @(\"93\u00f75=\", \"61\u00f77=\"),\n    @(\"35\u00f78=\", \"30\u00f74=\"),\n    @(\"21\u00f75=\", \"32\u00f73=\"),\n    @(\"29\u00f77=\", \"56\u00f79=\"),\n    @(\"80\u00f79=\", \"89\u00f78=\"),\n    @(\"40\u00f78=\", \"26\u00f79=\"),\n    @(\"18\u00f77=\", \"60\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
